$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.07419882731820041
$ws.Range("C2").Value = 1.337757721272759
$ws.Range("D2").Value = 9.827596144312645
$ws.Range("E2").Value = 3.134899702432702
$ws.Range("F2").Value = 3.165206548063491
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = -0.0318719375149274
$ws.Range("C3").Value = 1.484978161934237
$ws.Range("D3").Value = 9.229147191442671
$ws.Range("E3").Value = 3.03795115027261
$ws.Range("F3").Value = 3.068625194235299
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = -0.07763877361277893
$ws.Range("C4").Value = 1.390091861491201
$ws.Range("D4").Value = 6.971496353914397
$ws.Range("E4").Value = 2.640359133510894
$ws.Range("F4").Value = 2.666567550876379
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = -0.0027173285822412
$ws.Range("C5").Value = 1.635040883269796
$ws.Range("D5").Value = 9.649822969126083
$ws.Range("E5").Value = 3.106416419143783
$ws.Range("F5").Value = 3.139288264952675
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = -0.1033523239345538
$ws.Range("C6").Value = 1.503069249315635
$ws.Range("D6").Value = 8.686182792925141
$ws.Range("E6").Value = 2.947233074075605
$ws.Range("F6").Value = 2.97726366649755
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = -0.00824197970892185
$ws.Range("C7").Value = 1.667562176238278
$ws.Range("D7").Value = 8.839740055682562
$ws.Range("E7").Value = 2.973170034774762
$ws.Range("F7").Value = 3.006012189562184
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = -0.1262422412353678
$ws.Range("C8").Value = 1.56762830693314
$ws.Range("D8").Value = 8.721767756028379
$ws.Range("E8").Value = 2.95326391574278
$ws.Range("F8").Value = 2.983905233098791
$ws.Range("G8").Value = 45
